# Auto-generated Excel COM-interop script applying market-price refresh
# to the Leve profit tables (H..N columns) across multiple job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 759.7143
$ws.Range("I2").Value = 763.8
$ws.Range("J2").Value = 749.5
$ws.Range("K2").Value = 763.8
$ws.Range("L2").Value = 749.5
$ws.Range("M2").Value = -650.8
$ws.Range("N2").Value = -975.5

# Row 17
$ws.Range("H17").Value = 2743.8572
$ws.Range("J17").Value = 3117.8333
$ws.Range("L17").Value = 9353.499899999999
$ws.Range("N17").Value = -9689.499899999999

# Row 29
$ws.Range("H29").Value = 159.5
$ws.Range("I29").Value = 159.5
$ws.Range("K29").Value = 478.5
$ws.Range("M29").Value = -197.5

# Row 40
$ws.Range("H40").Value = 8949.833000000001
$ws.Range("J40").Value = 11266.333
$ws.Range("L40").Value = 11266.333
$ws.Range("N40").Value = -11616.333

# Row 43
$ws.Range("H43").Value = 4330.36
$ws.Range("I43").Value = 4246.125
$ws.Range("J43").Value = 4370
$ws.Range("K43").Value = 4246.125
$ws.Range("L43").Value = 4370
$ws.Range("M43").Value = -4177.125
$ws.Range("N43").Value = -4508

# Row 70
$ws.Range("H70").Value = 2783.1667
$ws.Range("I70").Value = 1625
$ws.Range("K70").Value = 4875
$ws.Range("M70").Value = -4605

# Row 73
$ws.Range("H73").Value = 2783.1667
$ws.Range("I73").Value = 1625
$ws.Range("K73").Value = 4875
$ws.Range("M73").Value = -3939

# Row 86
$ws.Range("H86").Value = 8269.143
$ws.Range("I86").Value = 7379.4
$ws.Range("J86").Value = 10493.5
$ws.Range("K86").Value = 7379.4
$ws.Range("L86").Value = 10493.5
$ws.Range("M86").Value = -6256.4
$ws.Range("N86").Value = -12739.5

# Row 89
$ws.Range("H89").Value = 8269.143
$ws.Range("I89").Value = 7379.4
$ws.Range("J89").Value = 10493.5
$ws.Range("K89").Value = 36897
$ws.Range("L89").Value = 52467.5
$ws.Range("M89").Value = -31281
$ws.Range("N89").Value = -63699.5

# Row 138
$ws.Range("H138").Value = 5960.75
$ws.Range("I138").Value = 3710
$ws.Range("J138").Value = 6523.4375
$ws.Range("K138").Value = 11130
$ws.Range("L138").Value = 19570.3125
$ws.Range("M138").Value = -5990
$ws.Range("N138").Value = -29850.3125

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 11606.793
$ws.Range("I2").Value = 13890.652
$ws.Range("K2").Value = 13890.652
$ws.Range("M2").Value = -13777.652

# Row 32
$ws.Range("H32").Value = 72261080
$ws.Range("I32").Value = 71470210
$ws.Range("J32").Value = 83333336
$ws.Range("K32").Value = 71470210
$ws.Range("L32").Value = 83333336
$ws.Range("M32").Value = -71469923
$ws.Range("N32").Value = -83333910

# Row 45
$ws.Range("H45").Value = 3703.7778
$ws.Range("I45").Value = 1183.75
$ws.Range("J45").Value = 5719.8
$ws.Range("K45").Value = 1183.75
$ws.Range("L45").Value = 5719.8
$ws.Range("M45").Value = -806.75
$ws.Range("N45").Value = -6473.8

# Row 61
$ws.Range("H61").Value = 3713.0815
$ws.Range("I61").Value = 2637.361
$ws.Range("K61").Value = 2637.361
$ws.Range("M61").Value = -2425.361

# Row 116
$ws.Range("H116").Value = 11606.793
$ws.Range("I116").Value = 13890.652
$ws.Range("K116").Value = 13890.652
$ws.Range("M116").Value = -11596.652

# Row 136
$ws.Range("H136").Value = 3713.0815
$ws.Range("I136").Value = 2637.361
$ws.Range("K136").Value = 7912.083
$ws.Range("M136").Value = -5362.083

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 11606.793
$ws.Range("I3").Value = 13890.652
$ws.Range("K3").Value = 13890.652
$ws.Range("M3").Value = -13776.652

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 20000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -19376
$ws.Range("N62").Value = -6248

# Row 65
$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 20000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 100000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -96880
$ws.Range("N65").Value = -31240

# Row 132
$ws.Range("H132").Value = 3611.754
$ws.Range("I132").Value = 2974.3845
$ws.Range("K132").Value = 8923.1535
$ws.Range("M132").Value = -6393.1535

# Row 141
$ws.Range("H141").Value = 151961
$ws.Range("J141").Value = 151961
$ws.Range("L141").Value = 151961
$ws.Range("N141").Value = -162321

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 314.5
$ws.Range("I2").Value = 87
$ws.Range("J2").Value = 360
$ws.Range("K2").Value = 522
$ws.Range("L2").Value = 2160
$ws.Range("M2").Value = -409
$ws.Range("N2").Value = -2386

# Row 38
$ws.Range("H38").Value = 881.10345
$ws.Range("I38").Value = 80.9375
$ws.Range("J38").Value = 1865.9231
$ws.Range("K38").Value = 242.8125
$ws.Range("L38").Value = 5597.7693
$ws.Range("M38").Value = 104.1875
$ws.Range("N38").Value = -6291.7693

# Row 40
$ws.Range("H40").Value = 1079.7273
$ws.Range("I40").Value = 234.625
$ws.Range("J40").Value = 3333.3333
$ws.Range("K40").Value = 938.5
$ws.Range("L40").Value = 13333.3332
$ws.Range("M40").Value = -869.5
$ws.Range("N40").Value = -13471.3332

# Row 57
$ws.Range("H57").Value = 399449
$ws.Range("J57").Value = 399449
$ws.Range("L57").Value = 1198347
$ws.Range("N57").Value = -1199465

# Row 92
$ws.Range("H92").Value = 491.33334
$ws.Range("I92").Value = 183
$ws.Range("J92").Value = 799.6667
$ws.Range("K92").Value = 549
$ws.Range("L92").Value = 2399.0001
$ws.Range("M92").Value = 699
$ws.Range("N92").Value = -4895.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# Row 38
$ws.Range("H38").Value = 25998
$ws.Range("J38").Value = 25998
$ws.Range("L38").Value = 25998
$ws.Range("N38").Value = -26924

# Row 40
$ws.Range("H40").Value = 7496.75
$ws.Range("J40").Value = 9993.5
$ws.Range("L40").Value = 9993.5
$ws.Range("N40").Value = -10295.5

# Row 44
$ws.Range("H44").Value = 20316.334
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 25474.5
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 25474.5
$ws.Range("N44").Value = -26666.5
$ws.Range("M44").Value = -9404

# Row 55
$ws.Range("H55").Value = 13898.8
$ws.Range("I55").Value = 6995
$ws.Range("J55").Value = 15624.75
$ws.Range("K55").Value = 6995
$ws.Range("L55").Value = 15624.75
$ws.Range("M55").Value = -6668
$ws.Range("N55").Value = -16278.75

# Row 80
$ws.Range("H80").Value = 167040370
$ws.Range("I80").Value = 560555.5
$ws.Range("K80").Value = 560555.5
$ws.Range("M80").Value = -559557.5

# Row 83
$ws.Range("H83").Value = 167040370
$ws.Range("I83").Value = 560555.5
$ws.Range("K83").Value = 2802777.5
$ws.Range("M83").Value = -2797785.5

# Row 97
$ws.Range("H97").Value = 83335100
$ws.Range("I97").Value = 2400
$ws.Range("J97").Value = 125001450
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 125001450
$ws.Range("N97").Value = -125002442
$ws.Range("M97").Value = -1904

# Row 113
$ws.Range("H113").Value = 9536.4
$ws.Range("I113").Value = 4562.2
$ws.Range("K113").Value = 4562.2
$ws.Range("M113").Value = -2392.2

# Row 126
$ws.Range("H126").Value = 4893.778
$ws.Range("I126").Value = 3650
$ws.Range("J126").Value = 5249.143
$ws.Range("K126").Value = 10950
$ws.Range("L126").Value = 15747.429
$ws.Range("M126").Value = -8480
$ws.Range("N126").Value = -20687.429

# Row 132
$ws.Range("H132").Value = 10329.549
$ws.Range("I132").Value = 8249.166999999999
$ws.Range("K132").Value = 24747.501
$ws.Range("M132").Value = -22217.501

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 494.7
$ws.Range("I16").Value = 494.7
$ws.Range("K16").Value = 494.7
$ws.Range("M16").Value = -324.7

# Row 46
$ws.Range("H46").Value = 4703.375
$ws.Range("I46").Value = 3639.5
$ws.Range("K46").Value = 3639.5
$ws.Range("M46").Value = -3451.5

# Row 132
$ws.Range("H132").Value = 6380.846
$ws.Range("I132").Value = 5586.654
$ws.Range("K132").Value = 16759.962
$ws.Range("M132").Value = -14229.962
